# TimeInPhase.xlsx correction ("Correción psp´s Clase Configuración")
#
# The source workbook was re-saved after the web-query sheet/connection was
# renamed from "excel(1)" to "excel" (the "(1)" suffix Excel appends when a
# file of that name already exists was cleaned up), the matching defined
# name was renamed to match, and the "Reporte generado a las ..." timestamp
# footer was refreshed to the new export time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet "excel(1)" -> "excel".
#    (Doing this first makes the defined name's RefersTo re-serialize as
#    bare `excel!$A$1:$E$27` instead of a quoted `'excel(1)'!...` sheet ref.)
$ws.Name = "excel"

# 2) Rename the workbook-scoped defined name "excel_1" -> "excel" so it
#    matches the sheet/query name everywhere.
$n = $wb.Names.Item(1)
$n.Name = "excel"

# 3) Refresh the "Reporte generado a las HH:MM AM/PM el D/M/YYYY" footer
#    text (row 25) to reflect the new export time, 01:33 PM instead of
#    11:22 AM, same date.
$ws.Range("A25").Value = "Reporte generado a las 01:33 PM el 5/12/2018"
